# Tarefas.xlsx - "Update e Delete quase prontos"
# Fills in the "Delete" task row (row 7) with its responsible person and
# status, clears out the leftover blank placeholder rows (10-13) and the
# now-unused leading cells of row 14, and leaves the selection where the
# author left off (D9:F9).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 ("Delete") was still blank in the Responsavel / Status columns.
$ws.Range("D7").Value = "Chaves"
$ws.Range("G7").Value = "Andamento (50%)"

# The trailing blank rows (10-13) under the table are no longer needed.
$ws.Range("A10:I13").Clear()

# Row 14 only keeps its formatted C14:D14 merged cell; A14/B14 are cleared.
$ws.Range("A14:B14").Clear()

# Leave the selection on D9:F9, matching where editing left off.
$ws.Range("D9:F9").Select() | Out-Null
